$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "B"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "C"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "D"

$ws.Range("B6").Select()
